$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (October -> November) ---
$ws.Range("A2").Value = "Year-to-Date through November 2016 and November 2015 (Thousand Tons)"
$ws.Range("B5").Value = "November 2016 YTD"
$ws.Range("C5").Value = "November 2015 YTD"
$ws.Range("E5").Value = "November 2016 YTD"
$ws.Range("F5").Value = "November 2015 YTD"
$ws.Range("G5").Value = "November 2016 YTD"
$ws.Range("H5").Value = "November 2015 YTD"
$ws.Range("I5").Value = "November 2016 YTD"
$ws.Range("J5").Value = "November 2015 YTD"
$ws.Range("K5").Value = "November 2016 YTD"
$ws.Range("L5").Value = "November 2015 YTD"

# --- Data updates ---
# Row 13
$ws.Range("B13").Value = 56
$ws.Range("C13").Value = 51
$ws.Range("D13").Value = 0.11
$ws.Range("K13").Value = 56
$ws.Range("L13").Value = 51

# Row 14
$ws.Range("B14").Value = 14
$ws.Range("D14").Value = 1.24
$ws.Range("K14").Value = 14

# Row 16
$ws.Range("B16").Value = 43
$ws.Range("C16").Value = 44
$ws.Range("D16").Value = -0.041
$ws.Range("K16").Value = 43
$ws.Range("L16").Value = 44

# Row 17
$ws.Range("B17").Value = 899
$ws.Range("C17").Value = 1180
$ws.Range("D17").Value = -0.24
$ws.Range("E17").Value = 461
$ws.Range("F17").Value = 656
$ws.Range("G17").Value = 391
$ws.Range("H17").Value = 475
$ws.Range("K17").Value = 46
$ws.Range("L17").Value = 48

# Row 19
$ws.Range("C19").Value = 348
$ws.Range("D19").Value = -0.54
$ws.Range("F19").Value = 348

# Row 20
$ws.Range("B20").Value = 302
$ws.Range("C20").Value = 317
$ws.Range("D20").Value = -0.049
$ws.Range("E20").Value = 264
$ws.Range("F20").Value = 280
$ws.Range("H20").Value = 13
$ws.Range("K20").Value = 36
$ws.Range("L20").Value = 25

# Row 21
$ws.Range("B21").Value = 392
$ws.Range("C21").Value = 464
$ws.Range("D21").Value = -0.16
$ws.Range("G21").Value = 390
$ws.Range("H21").Value = 463

# Row 22
$ws.Range("B22").Value = 46
$ws.Range("C22").Value = 50
$ws.Range("D22").Value = -0.084
$ws.Range("E22").Value = 38
$ws.Range("F22").Value = 28
$ws.Range("K22").Value = 8
$ws.Range("L22").Value = 22

# Row 23
$ws.Range("B23").Value = 24
$ws.Range("C23").Value = 17
$ws.Range("D23").Value = 0.43
$ws.Range("K23").Value = 23
$ws.Range("L23").Value = 15

# Row 24
$ws.Range("B24").Value = 24
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 0.43
$ws.Range("K24").Value = 23
$ws.Range("L24").Value = 15

# Row 31
$ws.Range("B31").Value = 727
$ws.Range("C31").Value = 555
$ws.Range("D31").Value = 0.31
$ws.Range("E31").Value = 703
$ws.Range("F31").Value = 527
$ws.Range("K31").Value = 24
$ws.Range("L31").Value = 28

# Row 34
$ws.Range("B34").Value = 703
$ws.Range("C34").Value = 527
$ws.Range("D34").Value = 0.33
$ws.Range("E34").Value = 703
$ws.Range("F34").Value = 527

# Row 35
$ws.Range("B35").Value = 24
$ws.Range("C35").Value = 28
$ws.Range("D35").Value = -0.14000000000000001
$ws.Range("K35").Value = 24
$ws.Range("L35").Value = 28

# Row 41
$ws.Range("B41").Value = 408
$ws.Range("C41").Value = 345
$ws.Range("D41").Value = 0.18
$ws.Range("E41").Value = 408
$ws.Range("F41").Value = 345

# Row 43
$ws.Range("B43").Value = 408
$ws.Range("C43").Value = 345
$ws.Range("D43").Value = 0.18
$ws.Range("E43").Value = 408
$ws.Range("F43").Value = 345

# Row 46
$ws.Range("B46").Value = 1671
$ws.Range("C46").Value = 1453
$ws.Range("D46").Value = 0.15
$ws.Range("E46").Value = 1598
$ws.Range("F46").Value = 1359
$ws.Range("K46").Value = 74
$ws.Range("L46").Value = 93

# Row 48
$ws.Range("B48").Value = 1631
$ws.Range("C48").Value = 1389
$ws.Range("D48").Value = 0.17
$ws.Range("E48").Value = 1598
$ws.Range("F48").Value = 1359
$ws.Range("K48").Value = 33
$ws.Range("L48").Value = 29

# Row 50
$ws.Range("B50").Value = 41
$ws.Range("C50").Value = 64
$ws.Range("D50").Value = -0.37
$ws.Range("K50").Value = 41
$ws.Range("L50").Value = 64

# Row 51
$ws.Range("B51").Value = 153
$ws.Range("C51").Value = 168
$ws.Range("D51").Value = -0.089
$ws.Range("G51").Value = 153
$ws.Range("H51").Value = 168

# Row 55
$ws.Range("B55").Value = 153
$ws.Range("C55").Value = 168
$ws.Range("D55").Value = -0.089
$ws.Range("G55").Value = 153
$ws.Range("H55").Value = 168

# Row 67
$ws.Range("B67").Value = 3939
$ws.Range("C67").Value = 3768
$ws.Range("D67").Value = 0.045
$ws.Range("E67").Value = 3170
$ws.Range("F67").Value = 2888
$ws.Range("G67").Value = 544
$ws.Range("H67").Value = 643
$ws.Range("K67").Value = 224
$ws.Range("L67").Value = 235
